# Apply updated loading_percent values (case with 380 kV) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 7.374395229308384
    "C2" = 5.203795761811899
    "E2" = 22.43784363511274
    "F2" = 38.4438280466111
    "G2" = 24.55215707027971
    "H2" = 13.36719475602441
    "I2" = 19.20001780696323
    "J2" = 7.802901303766303
    "K2" = 8.010754545515489
    "M2" = 17.7238058941527
    "N2" = 18.09932247087093
    "O2" = 19.76198620739943
    "B3" = 7.077448063182628
    "C3" = 5.096704019498737
    "E3" = 22.33542145644505
    "F3" = 38.38599728069765
    "G3" = 24.64510191976469
    "H3" = 13.41017223664875
    "I3" = 19.28218386667823
    "J3" = 7.811596499272455
    "K3" = 7.779556212349416
    "M3" = 17.59941163622951
    "N3" = 18.14824180939638
    "O3" = 19.83652868336744
    "B4" = 6.889245488179652
    "C4" = 5.029201761364854
    "E4" = 22.27712177456487
    "F4" = 38.35989117209337
    "G4" = 24.70940765283616
    "H4" = 13.43832495559749
    "I4" = 19.33578196842974
    "J4" = 7.817233551842567
    "K4" = 7.632630657485709
    "M4" = 17.52562396554414
    "N4" = 18.17979661808574
    "O4" = 19.88591780597908
    "B5" = 6.811185048778993
    "C5" = 5.001277900038029
    "E5" = 22.25453599160877
    "F5" = 38.35162517872366
    "G5" = 24.73742536140879
    "H5" = 13.45024152740713
    "I5" = 19.35841576266009
    "J5" = 7.819605865296452
    "K5" = 7.571562179373986
    "M5" = 17.49623206897485
    "N5" = 18.19303819089714
    "O5" = 19.90695388738918
    "B6" = 6.798144157934363
    "C6" = 4.99661673355544
    "E6" = 22.25085691285107
    "F6" = 38.35039610618799
    "G6" = 24.74218693915891
    "H6" = 13.45224710148193
    "I6" = 19.36222194768155
    "J6" = 7.820004331802015
    "K6" = 7.561351303296399
    "M6" = 17.49139319545729
    "N6" = 18.19526008909742
    "O6" = 19.91050183366237
    "B7" = 6.888198113476665
    "C7" = 5.028826823727199
    "E7" = 22.27681240741706
    "F7" = 38.35977007876599
    "G7" = 24.70977818012251
    "H7" = 13.43848386787178
    "I7" = 19.33608400757666
    "J7" = 7.817265241079976
    "K7" = 7.631811830384254
    "M7" = 17.5252248011903
    "N7" = 18.1799736475906
    "O7" = 19.88619782373262
    "B8" = 7.273290002085167
    "C8" = 5.167246813987339
    "E8" = 22.40158758889747
    "F8" = 38.42194209853333
    "G8" = 24.58269809676918
    "H8" = 13.3816474629691
    "I8" = 19.22769568973256
    "J8" = 7.805837638475041
    "K8" = 7.932099780906807
    "M8" = 17.68039273775248
    "N8" = 18.11587536308782
    "O8" = 19.78693681654588
    "B9" = 7.977599012635157
    "C9" = 5.423862795065766
    "E9" = 22.68180733605704
    "F9" = 38.61801363421952
    "G9" = 24.39122974276666
    "H9" = 13.28417084122097
    "I9" = 19.04009600582378
    "J9" = 7.785785179004176
    "K9" = 8.479389878658514
    "M9" = 18.00405863141747
    "N9" = 18.00217923592677
    "O9" = 19.62103398202681
    "B10" = 8.459218758925141
    "C10" = 5.602166956147967
    "E10" = 22.90809023444778
    "F10" = 38.80653941706633
    "G10" = 24.28615990173172
    "H10" = 13.22104910937764
    "I10" = 18.91743866181478
    "J10" = 7.772477214856393
    "K10" = 8.853584832421355
    "M10" = 18.25203763454151
    "N10" = 17.9258984410666
    "O10" = 19.51670678536665
    "B11" = 8.669731376000261
    "C11" = 5.680833858929668
    "E11" = 23.01516727104783
    "F11" = 38.90176201410328
    "G11" = 24.2461705683763
    "H11" = 13.19417214255724
    "I11" = 18.86492571586032
    "J11" = 7.766729789246961
    "K11" = 9.017303629490446
    "M11" = 18.36669832934271
    "N11" = 17.89275743311297
    "O11" = 19.47306693761943
    "B12" = 8.748158492768932
    "C12" = 5.710254223592901
    "E12" = 23.05628311356137
    "F12" = 38.93916071085518
    "G12" = 24.2321561882069
    "H12" = 13.18425836454625
    "I12" = 18.84551231258391
    "J12" = 7.764597256758055
    "K12" = 9.078331647383663
    "M12" = 18.41035272381431
    "N12" = 17.88043106010581
    "O12" = 19.45709158391645
    "B13" = 8.73132596184616
    "C13" = 5.70393470254143
    "E13" = 23.04740325903516
    "F13" = 38.93104698519058
    "G13" = 24.23512414519354
    "H13" = 13.1863817391831
    "I13" = 18.84967234188256
    "J13" = 7.765054586362181
    "K13" = 9.065231749836599
    "M13" = 18.40094103009161
    "N13" = 17.88307584233671
    "O13" = 19.46050768116724
    "B14" = 8.676209755600853
    "C14" = 5.683261782985327
    "E14" = 23.01853867015945
    "F14" = 38.90481209720507
    "G14" = 24.24499494303347
    "H14" = 13.19335124172799
    "I14" = 18.86331910270899
    "J14" = 7.766553465817258
    "K14" = 9.022344034568345
    "M14" = 18.37028527784571
    "N14" = 17.89173886237967
    "O14" = 19.47174160534036
    "B15" = 8.642280045468521
    "C15" = 5.670550458277464
    "E15" = 23.00093143366358
    "F15" = 38.88891633025625
    "G15" = 24.25118825385846
    "H15" = 13.19765462892071
    "I15" = 18.87173961513846
    "J15" = 7.767477283829092
    "K15" = 8.995946983389297
    "M15" = 18.35153738709682
    "N15" = 17.89707428159817
    "O15" = 19.47869437623633
    "B16" = 8.445283678502337
    "C16" = 5.596975167443835
    "E16" = 22.9011734997365
    "F16" = 38.80050506969258
    "G16" = 24.28893092894859
    "H16" = 13.22284253241787
    "I16" = 18.92093658114808
    "J16" = 7.772858974136931
    "K16" = 8.842751598218964
    "M16" = 18.24457888890606
    "N16" = 17.92809559238584
    "O16" = 19.51963568151249
    "B17" = 8.322193215840779
    "C17" = 5.551200091151691
    "E17" = 22.84101651879454
    "F17" = 38.74867689228675
    "G17" = 24.31408893621573
    "H17" = 13.23876490580546
    "I17" = 18.95195841314103
    "J17" = 7.776238829088235
    "K17" = 8.747080304646994
    "M17" = 18.17941566464308
    "N17" = 17.94752498176876
    "O17" = 19.5457306716591
    "B18" = 8.250590157480312
    "C18" = 5.524642351165885
    "E18" = 22.8068071712477
    "F18" = 38.7197585934527
    "G18" = 24.32929370091194
    "H18" = 13.24809599664057
    "I18" = 18.97011051215742
    "J18" = 7.778211685472583
    "K18" = 8.691442238244136
    "M18" = 18.14211191054511
    "N18" = 17.95884709079955
    "O18" = 19.56109923775536
    "B19" = 8.226210225090618
    "C19" = 5.515611560453584
    "E19" = 22.79529250180252
    "F19" = 38.71012113242946
    "G19" = 24.33456774114127
    "H19" = 13.25128505872669
    "I19" = 18.97630960639642
    "J19" = 7.778884621577367
    "K19" = 8.672500375602771
    "M19" = 18.12951275962344
    "N19" = 17.96270580802747
    "O19" = 19.56636447238667
    "B20" = 8.335380117566924
    "C20" = 5.556096771272833
    "E20" = 22.84738002195394
    "F20" = 38.75410192473194
    "G20" = 24.31133475582492
    "H20" = 13.23705204299186
    "I20" = 18.94862409061891
    "J20" = 7.775876052810557
    "K20" = 8.757328112179771
    "M20" = 18.18633437211844
    "N20" = 17.94544149908294
    "O20" = 19.54291561029702
    "B21" = 8.692434121506039
    "C21" = 5.689344070902812
    "E21" = 23.02700170355699
    "F21" = 38.9124817270172
    "G21" = 24.24206497191163
    "H21" = 13.19129696858993
    "I21" = 18.85929790672319
    "J21" = 7.76611201871157
    "K21" = 9.034967735546426
    "M21" = 18.37928349479107
    "N21" = 17.88918826687699
    "O21" = 19.46842699230577
    "B22" = 8.918253307978196
    "C22" = 5.774270346633389
    "E22" = 23.1476928307043
    "F22" = 39.02379152362545
    "G22" = 24.20337461966144
    "H22" = 13.16293174318303
    "I22" = 18.8036698321136
    "J22" = 7.759986420381278
    "K22" = 9.210761769818721
    "M22" = 18.50674108943005
    "N22" = 17.85372533314269
    "O22" = 19.42295126312638
    "B23" = 8.79843532689774
    "C23" = 5.729146567380436
    "E23" = 23.08298518243587
    "F23" = 38.96367703589026
    "G23" = 24.22342032077036
    "H23" = 13.17793012506325
    "I23" = 18.83310785293495
    "J23" = 7.7632324216556
    "K23" = 9.117465134809862
    "M23" = 18.43860119707929
    "N23" = 17.87253372977797
    "O23" = 19.4469287793723
    "B24" = 8.329420921141761
    "C24" = 5.55388373164486
    "E24" = 22.84450190982348
    "F24" = 38.75164652982779
    "G24" = 24.31257761279629
    "H24" = 13.23782587657161
    "I24" = 18.95013054964967
    "J24" = 7.776039971338422
    "K24" = 8.752697054862926
    "M24" = 18.183205925368
    "N24" = 17.94638296839515
    "O24" = 19.54418715888988
    "B25" = 7.793051357363262
    "C25" = 5.356156811711306
    "E25" = 22.60232048543361
    "F25" = 38.55709993718598
    "G25" = 24.43680277131183
    "H25" = 13.30904718586735
    "I25" = 19.08817970666092
    "J25" = 7.790958844249453
    "K25" = 8.336064921168415
    "M25" = 17.91458911873551
    "N25" = 18.03165899052741
    "O25" = 19.66283373401462
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
